$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking D-column price cells
# so Excel COM does not silently convert them to Number type
$textForceRefs = @("D5", "D6", "D9", "D15", "D17", "D19", "D20", "D21", "D25", "D27", "D29", "D31", "D32", "D36", "D37", "D44", "D46", "D47", "D48", "D49")
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "72.326.46"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "2.638.03"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "583.61"
$ws.Range("E5").Value = "  -3.09%  "
$ws.Range("D6").Value = "175.21"
$ws.Range("E6").Value = "  -1.56%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -1.06%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.172"
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("B10").Value = "LidoStakedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D10").Value = "2.638.43"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("E13").Value = "  -2.60%  "
$ws.Range("D14").Value = "3.122.68"
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("D15").Value = "0.0000185"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "72.247.87"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "25.82"
$ws.Range("E17").Value = "  -2.16%  "
$ws.Range("D18").Value = "2.650.71"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").Value = "12.09"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("D20").Value = "7.91"
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("D21").Value = "374.81"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "70.80"
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("E26").Value = "  -2.36%  "
$ws.Range("D27").Value = "9.50"
$ws.Range("E27").Value = "  -4.22%  "
$ws.Range("D28").Value = "2.777.25"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "0.0₃0947"
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").Value = "7.95"
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("D32").Value = "495.14"
$ws.Range("E32").Value = "  -4.31%  "
$ws.Range("E33").Value = "  -2.94%  "
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "162.67"
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("D37").Value = "19.16"
$ws.Range("E37").Value = "  -1.83%  "
$ws.Range("E38").Value = "  +4.85%  "
$ws.Range("E39").Value = "  -1.46%  "
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("E42").Value = "  -6.18%  "
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("D44").Value = "4.89"
$ws.Range("E44").Value = "  -2.96%  "
$ws.Range("E45").Value = "  -2.49%  "
$ws.Range("D46").Value = "39.02"
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("D47").Value = "151.37"
$ws.Range("E47").Value = "  -1.91%  "
$ws.Range("D48").Value = "3.64"
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("D49").Value = "0.543"
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("E50").Value = "  -3.72%  "
$ws.Range("E51").Value = "  -0.71%  "

# Restore default number format on the cells we temporarily forced to text
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).NumberFormat = "General"
}
